$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("r0")

$row = 34
$ws.Cells.Item($row, 1).Value = "ExoT_r0_script_13v0"
$ws.Cells.Item($row, 2).Value = "Griffin qiazol elution"
$ws.Cells.Item($row, 3).Value = "20 mL"
$ws.Cells.Item($row, 4).Value = "5 mL"
$ws.Cells.Item($row, 5).Value = "1 hour"
$ws.Cells.Item($row, 6).Value = "2.5 mL"
$ws.Cells.Item($row, 7).Value = "15 mL/hr"
$ws.Cells.Item($row, 8).Value = "15 mL/hr"
$ws.Cells.Item($row, 9).Value = "200-800-1000"
$ws.Cells.Item($row, 10).Value = "2 mins"
$ws.Cells.Item($row, 11).Value = "N"

$wb.Save()
